$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.660.34"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "2.302.65"
$ws.Range("E3").Value = "  +0.15%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'316.19"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "'103.80"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.629"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  +0.16%  "

$ws.Range("D9").Value = "'0.606"
$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").Value = "'39.87"
$ws.Range("E10").Value = "  +1.15%  "

$ws.Range("D11").Value = "'0.0905"
$ws.Range("E11").Value = "  -0.29%  "

$ws.Range("D12").Value = "'8.53"
$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("E13").Value = "  +0.82%  "

$ws.Range("D14").Value = "'0.998"
$ws.Range("E14").Value = "  +3.87%  "

$ws.Range("D15").Value = "'15.35"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").Value = "2.652.62"
$ws.Range("E16").Value = "  +0.23%  "

$ws.Range("D17").Value = "2.309.56"
$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("D18").Value = "42.605.21"
$ws.Range("E18").Value = "  +1.02%  "

$ws.Range("D19").Value = "'7.60"
$ws.Range("E19").Value = "  +2.92%  "

$ws.Range("E20").Value = "  +0.36%  "

$ws.Range("D21").Value = "'13.65"
$ws.Range("E21").Value = "  +31.55%  "

$ws.Range("D22").Value = "'73.98"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("E23").Value = "  -1.51%  "

$ws.Range("D24").Value = "'268.22"
$ws.Range("E24").Value = "  -4.17%  "

$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").Value = "'10.90"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").Value = "'22.56"
$ws.Range("E29").Value = "  -1.31%  "

$ws.Range("D30").Value = "'6.67"
$ws.Range("E30").Value = "  +14.85%  "

$ws.Range("D31").Value = "'37.76"
$ws.Range("E31").Value = "  +4.35%  "

$ws.Range("D32").Value = "'165.77"
$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("D33").Value = "'0.0882"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").Value = "'2.67"
$ws.Range("E34").Value = "  -5.61%  "

$ws.Range("D36").Value = "'0.114"
$ws.Range("E36").Value = "  +0.43%  "

$ws.Range("D37").Value = "'4.58"
$ws.Range("E37").Value = "  +1.50%  "

$ws.Range("E38").Value = "  +1.88%  "

$ws.Range("D39").Value = "'3.75"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("D40").Value = "'2.73"
$ws.Range("E40").Value = "  -3.10%  "

$ws.Range("D41").Value = "'1.62"
$ws.Range("E41").Value = "  +11.55%  "

$ws.Range("D42").Value = "'98.45"
$ws.Range("E42").Value = "  -1.61%  "

$ws.Range("D43").Value = "'70.05"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("D44").Value = "'0.226"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  +0.02%  "

$ws.Range("E46").Value = "  +3.70%  "

$ws.Range("D47").Value = "'116.76"
$ws.Range("E47").Value = "  +3.90%  "

$ws.Range("D48").Value = "'80.67"
$ws.Range("E48").Value = "  +4.35%  "

$ws.Range("D49").Value = "1.635.86"
$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("D51").Value = "'8.89"
$ws.Range("E51").Value = "  -0.18%  "
